$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new "Preparation" header column (H) to the first worksheet,
# mirroring the "Prep Needed?" column already present on sheet 2.
$ws1.Range("H2").Value = "Preparation"

# Copy the formatting from the neighboring header cell (G2) so the new
# header cell (H2) matches the bold/shaded header style used by the rest
# of the header row.
$ws1.Range("G2").Copy()
$ws1.Range("H2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to reflect where the user left off editing.
$ws1.Activate()
$ws1.Range("H7").Select()
